$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("M2").Value = 13

$ws.Range("H4").Value = 2105.818
$ws.Range("I4").Value = 2105.818
$ws.Range("K4").Value = 2105.818
$ws.Range("M4").Value = -1991.818

$ws.Range("H5").Value = 104.5
$ws.Range("I5").Value = 87.35714
$ws.Range("K5").Value = 87.35714
$ws.Range("M5").Value = 27.64286

$ws.Range("H15").Value = 2084.0833
$ws.Range("I15").Value = 2084.0833
$ws.Range("K15").Value = 6252.249899999999
$ws.Range("M15").Value = -6083.249899999999

$ws.Range("H97").Value = 4282.5713
$ws.Range("J97").Value = 4282.5713
$ws.Range("L97").Value = 12847.7139
$ws.Range("N97").Value = -13839.7139

$ws.Range("H116").Value = 4983
$ws.Range("I116").Value = 4987
$ws.Range("J116").Value = 4975
$ws.Range("K116").Value = 4987
$ws.Range("L116").Value = 4975
$ws.Range("M116").Value = -1545
$ws.Range("N116").Value = -11859

$ws.Range("H127").Value = 796
$ws.Range("I127").Value = 796
$ws.Range("K127").Value = 2388
$ws.Range("M127").Value = 2572

$ws.Range("H137").Value = 1501.3334
$ws.Range("I137").Value = 1258.3572
$ws.Range("J137").Value = 4903
$ws.Range("K137").Value = 3775.0716
$ws.Range("L137").Value = 14709
$ws.Range("M137").Value = -1225.0716
$ws.Range("N137").Value = -19809

$ws.Range("H138").Value = 3499.8333
$ws.Range("I138").Value = 2500
$ws.Range("J138").Value = 3999.75
$ws.Range("K138").Value = 7500
$ws.Range("L138").Value = 11999.25
$ws.Range("M138").Value = -2360
$ws.Range("N138").Value = -22279.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H74").Value = 5659.8
$ws.Range("I74").Value = 5659.8
$ws.Range("K74").Value = 5659.8
$ws.Range("M74").Value = -4785.8

$ws.Range("H77").Value = 5659.8
$ws.Range("I77").Value = 5659.8
$ws.Range("K77").Value = 28299
$ws.Range("M77").Value = -23931

$ws.Range("H112").Value = 60000
$ws.Range("J112").Value = 60000
$ws.Range("L112").Value = 60000
$ws.Range("N112").Value = -62954

$ws.Range("H122").Value = 2535.375
$ws.Range("I122").Value = 2183.2856
$ws.Range("K122").Value = 6549.8568
$ws.Range("M122").Value = -4099.8568

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1963.9
$ws.Range("J20").Value = 1249.5
$ws.Range("L20").Value = 1249.5
$ws.Range("N20").Value = -1743.5

$ws.Range("H22").Value = 349.75
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H7").Value = 142.58333
$ws.Range("I7").Value = 127.75
$ws.Range("J7").Value = 172.25
$ws.Range("K7").Value = 127.75
$ws.Range("L7").Value = 172.25
$ws.Range("M7").Value = -14.75
$ws.Range("N7").Value = -398.25

$ws.Range("H31").Value = 17803.264
$ws.Range("I31").Value = 10591
$ws.Range("J31").Value = 37997.6
$ws.Range("K31").Value = 10591
$ws.Range("L31").Value = 37997.6
$ws.Range("M31").Value = -10296
$ws.Range("N31").Value = -38587.6

$ws.Range("H34").Value = 17803.264
$ws.Range("I34").Value = 10591
$ws.Range("J34").Value = 37997.6
$ws.Range("K34").Value = 10591
$ws.Range("L34").Value = 37997.6
$ws.Range("M34").Value = -10389
$ws.Range("N34").Value = -38401.6

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H59").Value = 35571.43

$ws.Range("H132").Value = 5084.625
$ws.Range("I132").Value = 5367.1665
$ws.Range("J132").Value = 4237
$ws.Range("K132").Value = 16101.4995
$ws.Range("L132").Value = 12711
$ws.Range("M132").Value = -13571.4995
$ws.Range("N132").Value = -17771

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H141").Value = 488772
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 488772
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 488772
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -499132

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 51.5
$ws.Range("J44").Value = 100
$ws.Range("L44").Value = 300
$ws.Range("N44").Value = -1096

$ws.Range("H70").Value = 6999
$ws.Range("I70").Value = 6999
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 20997
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -20682
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 6999
$ws.Range("I73").Value = 6999
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 20997
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -19905
$ws.Range("N73").ClearContents()

$ws.Range("H109").Value = 506
$ws.Range("I109").Value = 506
$ws.Range("K109").Value = 1518
$ws.Range("M109").Value = -478

$ws.Range("H121").Value = 1383.5
$ws.Range("I121").Value = 300.4
$ws.Range("K121").Value = 901.1999999999999
$ws.Range("M121").Value = 408.8000000000001

$ws.Range("H122").Value = 2027.4615
$ws.Range("J122").Value = 2313.625
$ws.Range("L122").Value = 20822.625
$ws.Range("N122").Value = -25722.625

$ws.Range("H131").Value = 1826.875
$ws.Range("J131").Value = 1830.8334
$ws.Range("L131").Value = 5492.5002
$ws.Range("N131").Value = -15572.5002

$ws.Range("H132").Value = 3083.1667
$ws.Range("I132").Value = 3125
$ws.Range("K132").Value = 28125
$ws.Range("M132").Value = -25595

$ws.Range("H139").Value = 1899.25
$ws.Range("I139").Value = 1199
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 3597
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = 1543
$ws.Range("N139").Value = -22280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 4167633.2
$ws.Range("I14").Value = 5000766.5
$ws.Range("J14").Value = 3334500
$ws.Range("K14").Value = 5000766.5
$ws.Range("L14").Value = 3334500
$ws.Range("M14").Value = -5000598.5
$ws.Range("N14").Value = -3334836

$ws.Range("H97").Value = 828.5454999999999
$ws.Range("I97").Value = 846.4
$ws.Range("K97").Value = 846.4
$ws.Range("M97").Value = -350.4

$ws.Range("H102").Value = 3238.5715
$ws.Range("I102").Value = 3238.5715
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3238.5715
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1616.5715
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 16557.5
$ws.Range("I122").Value = 1123.75
$ws.Range("K122").Value = 3371.25
$ws.Range("M122").Value = -921.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3220.2
$ws.Range("I61").Value = 3824
$ws.Range("K61").Value = 3824
$ws.Range("M61").Value = -3622

$ws.Range("H82").Value = 1868.5
$ws.Range("I82").Value = 1999.7142
$ws.Range("K82").Value = 1999.7142
$ws.Range("M82").Value = -1638.7142

$ws.Range("H85").Value = 1868.5
$ws.Range("I85").Value = 1999.7142
$ws.Range("K85").Value = 1999.7142
$ws.Range("M85").Value = -751.7141999999999

$ws.Range("H113").Value = 3220.2
$ws.Range("I113").Value = 3824
$ws.Range("K113").Value = 3824
$ws.Range("M113").Value = -1654

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H126").Value = 2050.75
$ws.Range("I126").Value = 2050.75
$ws.Range("K126").Value = 6152.25
$ws.Range("M126").Value = -3682.25
